$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the monster name in B4 from "火焰哥布林" (Flame Goblin) to "寒冰哥布林" (Ice Goblin)
$ws.Range("B4").Value = "寒冰哥布林"

# Reflect the user's final selection after the edit
$ws.Range("B4").Select()
